$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 230.71428
$ws.Range("I2").Value = 222.5
$ws.Range("J2").Value = 241.66667
$ws.Range("K2").Value = 222.5
$ws.Range("L2").Value = 241.66667
$ws.Range("M2").Value = -109.5
$ws.Range("N2").Value = -467.66667
$ws.Range("H19").Value = 821.8461
$ws.Range("I19").Value = 405.7143
$ws.Range("J19").Value = 1307.3334
$ws.Range("K19").Value = 405.7143
$ws.Range("L19").Value = 1307.3334
$ws.Range("M19").Value = -230.7143
$ws.Range("N19").Value = -1657.3334
$ws.Range("H40").Value = 2255.6155
$ws.Range("I40").Value = 1564.6
$ws.Range("J40").Value = 2687.5
$ws.Range("K40").Value = 1564.6
$ws.Range("L40").Value = 2687.5
$ws.Range("M40").Value = -1389.6
$ws.Range("N40").Value = -3037.5
$ws.Range("H74").Value = 3442.5173
$ws.Range("I74").Value = 3125
$ws.Range("J74").Value = 3493.32
$ws.Range("K74").Value = 3125
$ws.Range("L74").Value = 3493.32
$ws.Range("M74").Value = -2189
$ws.Range("N74").Value = -5365.32
$ws.Range("H76").Value = 3615.25
$ws.Range("I76").Value = 2999
$ws.Range("K76").Value = 2999
$ws.Range("M76").Value = -2684
$ws.Range("H77").Value = 3442.5173
$ws.Range("I77").Value = 3125
$ws.Range("J77").Value = 3493.32
$ws.Range("K77").Value = 15625
$ws.Range("L77").Value = 17466.6
$ws.Range("M77").Value = -10945
$ws.Range("N77").Value = -26826.6
$ws.Range("H79").Value = 3615.25
$ws.Range("I79").Value = 2999
$ws.Range("K79").Value = 2999
$ws.Range("M79").Value = -1907
$ws.Range("H129").Value = 800.39
$ws.Range("I129").Value = 401.75
$ws.Range("K129").Value = 1205.25
$ws.Range("M129").Value = 3794.75

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2600
$ws.Range("I63").Value = 2700
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 2700
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -2014
$ws.Range("N63").Value = -3872
$ws.Range("H66").Value = 2600
$ws.Range("I66").Value = 2700
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 13500
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -10068
$ws.Range("N66").Value = -19364

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 466.66666
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -427
$ws.Range("N22").Value = -746
$ws.Range("H86").Value = 2021.5938
$ws.Range("I86").Value = 1970.3182
$ws.Range("J86").Value = 2134.4
$ws.Range("K86").Value = 1970.3182
$ws.Range("L86").Value = 2134.4
$ws.Range("M86").Value = -847.3181999999999
$ws.Range("N86").Value = -4380.4
$ws.Range("H89").Value = 2021.5938
$ws.Range("I89").Value = 1970.3182
$ws.Range("J89").Value = 2134.4
$ws.Range("K89").Value = 9851.591
$ws.Range("L89").Value = 10672
$ws.Range("M89").Value = -4235.591
$ws.Range("N89").Value = -21904
$ws.Range("H105").Value = 2803.7856
$ws.Range("I105").Value = 1522.375
$ws.Range("J105").Value = 3017.3542
$ws.Range("K105").Value = 1522.375
$ws.Range("L105").Value = 3017.3542
$ws.Range("M105").Value = 224.625
$ws.Range("N105").Value = -6511.3542
$ws.Range("H133").Value = 30700
$ws.Range("J133").Value = 30700
$ws.Range("L133").Value = 30700
$ws.Range("N133").Value = -40820

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3211.2
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 3528
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 3528
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4776
$ws.Range("H65").Value = 3211.2
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 3528
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 17640
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -23880
$ws.Range("H105").Value = 785.1667
$ws.Range("I105").Value = 740
$ws.Range("K105").Value = 740
$ws.Range("M105").Value = 1007

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 513.75
$ws.Range("I110").Value = 513.75
$ws.Range("K110").Value = 1541.25
$ws.Range("M110").Value = 2548.75
$ws.Range("H131").Value = 708.3
$ws.Range("I131").Value = 353
$ws.Range("J131").Value = 882.1702
$ws.Range("K131").Value = 1059
$ws.Range("L131").Value = 2646.5106
$ws.Range("M131").Value = 3981
$ws.Range("N131").Value = -12726.5106

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6553.72
$ws.Range("I70").Value = 7092.2
$ws.Range("J70").Value = 4399.8
$ws.Range("K70").Value = 7092.2
$ws.Range("L70").Value = 4399.8
$ws.Range("M70").Value = -6822.2
$ws.Range("N70").Value = -4939.8
$ws.Range("H73").Value = 6553.72
$ws.Range("I73").Value = 7092.2
$ws.Range("J73").Value = 4399.8
$ws.Range("K73").Value = 7092.2
$ws.Range("L73").Value = 4399.8
$ws.Range("M73").Value = -6156.2
$ws.Range("N73").Value = -6271.8
$ws.Range("H80").Value = 2889.1714
$ws.Range("I80").Value = 2738.524
$ws.Range("J80").Value = 3115.1428
$ws.Range("K80").Value = 2738.524
$ws.Range("L80").Value = 3115.1428
$ws.Range("M80").Value = -1740.524
$ws.Range("N80").Value = -5111.1428
$ws.Range("H83").Value = 2889.1714
$ws.Range("I83").Value = 2738.524
$ws.Range("J83").Value = 3115.1428
$ws.Range("K83").Value = 13692.62
$ws.Range("L83").Value = 15575.714
$ws.Range("M83").Value = -8700.619999999999
$ws.Range("N83").Value = -25559.714

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3333.3823
$ws.Range("I122").Value = 3055.9614
$ws.Range("J122").Value = 4235
$ws.Range("K122").Value = 9167.8842
$ws.Range("L122").Value = 12705
$ws.Range("M122").Value = -6717.8842
$ws.Range("N122").Value = -17605
$ws.Range("H132").Value = 2255.782
$ws.Range("I132").Value = 1939.3334
$ws.Range("J132").Value = 3996.25
$ws.Range("K132").Value = 5818.0002
$ws.Range("L132").Value = 11988.75
$ws.Range("M132").Value = -3288.0002
$ws.Range("N132").Value = -17048.75
